$d = $word.ActiveDocument

# The paragraph we need to rewrite ("Projet devops" -> "Projet devops
# «  securit » ", with the word "devops" flagged by the spell checker)
# happens to be the last paragraph in the document body. InsertXML can only
# cleanly "split" a paragraph in place when there is a following paragraph
# to absorb the trailing paragraph mark; otherwise it leaves a stray empty
# paragraph behind. So we temporarily add a throw-away paragraph after it,
# do the rewrite, then remove the scratch paragraph again.
$d.Content.InsertParagraphAfter()

# Locate the (now not-last) paragraph that still reads "Projet devops".
$found = $d.Content.Duplicate
$found.Find.Execute("Projet devops", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$rewriteRange = $found.Paragraphs(1).Range
$paraIndex = $rewriteRange.Paragraphs(1).Index

# Full replacement markup for that paragraph: "Projet " + "devops" (marked
# as a spell-check error, like Word's proofer would do for a dictionary
# miss) + " «  securit » ", keeping the existing _GoBack bookmark.
$openXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="001F36BB" w:rsidRDefault="001F36BB" w:rsidP="001F36BB">
            <w:r><w:t xml:space="preserve">Projet </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>devops</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> «  securit » </w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$rewriteRange.InsertXML($openXml)

# Drop the scratch paragraph introduced above: delete the span from the end
# of the rewritten paragraph through the end of the now-redundant empty
# paragraph that trails it (i.e. its paragraph mark).
$rewrittenEnd = $d.Paragraphs($paraIndex).Range.End
$scratchEnd = $d.Paragraphs($paraIndex + 1).Range.End
$cleanup = $d.Range($rewrittenEnd, $scratchEnd)
$cleanup.Delete()
